$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wifi configuration")
$fcs = $ws.Range("E1").FormatConditions
for ($i=1; $i -le $fcs.Count; $i++) {
  $fc = $fcs.Item($i)
  $rng = $ws.Range("E1,E14:E1048576")
  Write-Host "rng addr" $rng.Address()
  $fc.ModifyAppliesToRange($rng)
}
